$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.159.37"
$ws.Range("E2").Value = "  -1.34%  "

# Row 3
$ws.Range("D3").Value = "2.765.77"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'354.42"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
$ws.Range("D6").Value = "'107.99"
$ws.Range("E6").Value = "  -0.85%  "

# Row 7
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  -2.31%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  -0.86%  "

# Row 10
$ws.Range("D10").Value = "'39.41"
$ws.Range("E10").Value = "  -1.47%  "

# Row 11
$ws.Range("E11").Value = "  +3.62%  "

# Row 12
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'20.05"
$ws.Range("E12").Value = "  +3.71%  "

# Row 13
$ws.Range("B13").Value = "Dogecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D13").Value = "'0.0833"
$ws.Range("E13").Value = "  -1.98%  "

# Row 14
$ws.Range("D14").Value = "'7.53"
$ws.Range("E14").Value = "  -0.61%  "

# Row 15
$ws.Range("D15").Value = "3.196.51"
$ws.Range("E15").Value = "  -0.41%  "

# Row 16
$ws.Range("D16").Value = "2.738.70"
$ws.Range("E16").Value = "  -2.22%  "

# Row 17
$ws.Range("D17").Value = "'0.930"
$ws.Range("E17").Value = "  +0.99%  "

# Row 18
$ws.Range("D18").Value = "51.133.19"
$ws.Range("E18").Value = "  -1.22%  "

# Row 19
$ws.Range("D19").Value = "'7.73"
$ws.Range("E19").Value = "  +5.37%  "

# Row 20
$ws.Range("D20").Value = "'3.07"
$ws.Range("E20").Value = "  -1.41%  "

# Row 21
$ws.Range("D21").Value = "'13.11"
$ws.Range("E21").Value = "  +0.92%  "

# Row 22
$ws.Range("E22").Value = "  -1.20%  "

# Row 23
$ws.Range("D23").Value = "'69.67"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("D24").Value = "'265.23"
$ws.Range("E24").Value = "  -2.73%  "

# Row 25
$ws.Range("D25").Value = "'2.72"
$ws.Range("E25").Value = "  -0.40%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").Value = "'25.96"
$ws.Range("E27").Value = "  -1.82%  "

# Row 28
$ws.Range("D28").Value = "'0.162"
$ws.Range("E28").Value = "  +13.03%  "

# Row 29
$ws.Range("E29").Value = "  +0.83%  "

# Row 30
$ws.Range("E30").Value = "  +0.94%  "

# Row 31
$ws.Range("D31").Value = "'51.82"
$ws.Range("E31").Value = "  +1.02%  "

# Row 32
$ws.Range("D32").Value = "'34.69"
$ws.Range("E32").Value = "  +2.78%  "

# Row 33
$ws.Range("D33").Value = "'6.04"
$ws.Range("E33").Value = "  +6.14%  "

# Row 34
$ws.Range("E34").Value = "  -3.40%  "

# Row 35
$ws.Range("D35").Value = "'5.47"
$ws.Range("E35").Value = "  +2.19%  "

# Row 36
$ws.Range("D36").Value = "'0.0830"
$ws.Range("E36").Value = "  -0.50%  "

# Row 37
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$ws.Range("D38").Value = "'18.19"
$ws.Range("E38").Value = "  +0.50%  "

# Row 39
$ws.Range("D39").Value = "'3.14"
$ws.Range("E39").Value = "  -1.61%  "

# Row 40
$ws.Range("E40").Value = "  -1.63%  "

# Row 41
$ws.Range("D41").Value = "'2.52"
$ws.Range("E41").Value = "  +0.48%  "

# Row 42
$ws.Range("E42").Value = "  -0.30%  "

# Row 43
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'22.22"
$ws.Range("E43").Value = "  +2.37%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'120.32"
$ws.Range("E44").Value = "  -2.44%  "

# Row 45
$ws.Range("E45").Value = "  -2.41%  "

# Row 46
$ws.Range("D46").Value = "2.085.36"
$ws.Range("E46").Value = "  +1.14%  "

# Row 47
$ws.Range("E47").Value = "  +0.13%  "

# Row 49
$ws.Range("E49").Value = "  -3.11%  "

# Row 50
$ws.Range("D50").Value = "'0.916"
$ws.Range("E50").Value = "  -0.42%  "

# Row 51
$ws.Range("E51").Value = "  +6.21%  "
